$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "Our results have demonstrated that ..." (Key results S1:S3):
#    a few small text edits.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("after 25 years", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "half way through the study period", 2) | Out-Null

$d.Content.Find.Execute("start of the simulation period", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "start of the study period", 2) | Out-Null

$d.Content.Find.Execute("The loss of trees during periods of low funding", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The rate of forest loss during periods of low funding", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Rework the "Conservation projects that are initially underfunded..."
#    sentence, then drop the empty paragraph + the two following paragraphs
#    ("Human populations are increasing..." / "Likewise, increasing manager
#    budgets are great...") whose content has been folded into the rewrite.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute( `
    "Conservation projects that are initially underfunded will spend many years working to reach the same levels of protection as they would have had, had they been well-funded at the start.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Conservation projects that are initially underfunded yet receive increasing resources will still spend many years working to reach the same levels of protection as they would have had, had they been provided an adequate, stable budget at the start. Our results suggest it could be several decades before the deforestation trajectories of the two alternative projects meet, and the increasing budget starts to have an effect.", `
    2) | Out-Null

# Delete the now-superfluous paragraphs (work backwards so indices stay put):
#   8 = "Likewise, increasing manager budgets are great..."
#   7 = "Human populations are increasing, as is pressure on landscapes..."
#   6 = (empty paragraph)
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(6).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Add two blank paragraphs before "These scenarios highlight ..."
#    (now paragraph 6 after the deletions above).
# ---------------------------------------------------------------------------

$d.Paragraphs.Item(6).Range.InsertParagraphBefore()
$d.Paragraphs.Item(6).Range.InsertParagraphBefore()

$d.Save()
